$d = $word.ActiveDocument

$d.Content.Find.Execute('Ativação: 01/01/2021', $true, $false, $false, $false, $false, $true, 1, $false, 'Ativação: 01/01/2024', 2) | Out-Null
$d.Content.Find.Execute('Conferir aos alunos uma visão geral da indústria química e correlatas, bem como das principais características dos processos e arranjos produtivos desta indústria.', $true, $false, $false, $false, $false, $true, 1, $false, 'Conferir aos alunos uma visão geral da indústria química e correlatas, processos e produtos, e uma visão global das matérias primas mais importantes da indústria química.', 2) | Out-Null
$d.Content.Find.Execute('Providing to the students an overview of the chemical and related industries, as well as the main features of the processes and production arrangements of this industry.', $true, $false, $false, $false, $false, $true, 1, $false, 'Objectives:Check the students an overview of the chemical industry and related industries, processes and products, and an overview of the most important raw materials in the chemical industry..', 2) | Out-Null
$d.Content.Find.Execute('Processos da Indústria Químicos.', $true, $false, $false, $false, $false, $true, 1, $false, 'Processo Químico e Indústria Química; Química Fina; Petróleo, Gás Natural e Petroquímica; Plásticos e afins; Fertilizantes; Vidro; Celulose e Papel.', 2) | Out-Null
$d.Content.Find.Execute('Current and relevant topics related to chemical processes.', $true, $false, $false, $false, $false, $true, 1, $false, 'Chemical Process and Chemical Industry; Fine Chemistry; Oil, Natural Gas and Petrochemicals; Plastics and allied products; fertilizers; Glass; Cellulose and paper.', 2) | Out-Null
$d.Content.Find.Execute('Panorama da Indústria Química. Química Fina. Petróleo e Petroquímica. Cerâmica. Vidro. Cimento. Celulose e Papel. Plásticos e afins. Processos Químicos Inovadores.', $true, $false, $false, $false, $false, $true, 1, $false, '1- Processo Químico e Indústria Química; 2- Química Fina: 2.1- Características, 2.2- Principais Segmentos (Defensivos Agrícolas, Fármacos, Catalisadores, Corantes e Pigmentos, Especialidades), 2.3- Química Fina X Química de Base, 3- Petróleo, Gás Natural e Petroquímica; 4- Plásticos e Afins: 4.1- Resinas e suas respectivas aplicações; 4.2- Reciclagem de plásticos; 5- Fertilizantes; 6- Vidro; 7- Celulose e papel.', 2) | Out-Null
$d.Content.Find.Execute('Overview of the Chemical Industry. Fine Chemistry. Petroleum and Petrochemicals. Ceramics. Glass. Cement. Cellulose and paper. Plastics and allied products. Innovative Chemical Processes.', $true, $false, $false, $false, $false, $true, 1, $false, '1- Fine Chemicals: 1.1- Characteristics, 1.2- Main Segments (Agricultural Defensives, Pharmaceuticals, Catalysts, Dyes and Pigments, Specialties), 1.3- Fine Chemicals X Basic Chemicals, 2- Oil, Natural Gas and Petrochemicals; 3- Unitary Processes Organic: 3.1- Alkylation and Acylation; 3.2- Hydrogenation and Dehydrogenation; 3.2.1 Oxo processes, 3.2.2- Amino; 3.3- Halogenation; 3.4- Esterification; 3.5- Sulfonation/Sulfation; 3.6- Oxidation.', 2) | Out-Null
$d.Content.Find.Execute('Provas e/ou trabalhos.', $true, $false, $false, $false, $false, $true, 1, $false, 'Provas em sala, apresentações em sala, entrega de exercícios ou casos práticos elaborados fora de sala de aula e frequência.', 2) | Out-Null
$d.Content.Find.Execute('Material elaborado pelo docente.   LIVROS:Ullmann’s encyclopedia of industrial chemistry; Editorial advisory board, Giuseppe Bellussi et al.; 7th, completely revised edition; Weinheim; New York: WileyVCH, 2011.  Encyclopedia of Chemical Processing; Edited by Sunggyu Lee; New York: Taylor & Francis, 2006.Kirk, Raymond Eller. Encyclopedia of chemical technology / Herman F.Mark et al. New York: John Wiley, 1984.   Shreve, R. Norris; BRINK JR., J. A. Indústrias de processos químicos. Tradução de Horácio Macedo; 4.ed. Rio de Janeiro: Editora Guanabara Koogan, 2008, c1997.   REVISTAS:Química & Derivados. Disponível em: http://www.quimica.com.br/pquimica/category/revista/Petróleo & Energia. Disponível em: http://www.petroleoenergia.com.br/petroleo/category/revista-petroleo-e-energia/.', $true, $false, $false, $false, $false, $true, 1, $false, 'Material elaborado pelo docente;Livros:Ullmann’s encyclopedia of industrial chemistry; Editorial advisory board, Giuseppe Bellussi et al.; 7th, completely revised edition; Weinheim ; New York : WileyVCH, 2011.Encyclopedia of Chemical Processing; Edited by Sunggyu Lee; New York : Taylor & Francis, 2006.Kirk, Raymond Eller. Encyclopedia of chemical technology / Herman F.Mark et al. New York: John Wiley, 1984.Manual Econômico da Indústria Química - MEIQ / Centro de Pesquisas e Desenvolvimento; 8ed; Camaçari: CEPED, 2007.Shreve, R. Norris; BRINK JR., J. A. Indústrias de processos químicos. Tradução de Horácio Macedo; 4.ed. Rio de Janeiro: Editora Guanabara Koogan, 2008, c1997.Revistas:Química & Derivados, São Paulo, SP: QD, v. 1, n. 1, nov. 1965-; Disponível em: http://www.quimica.com.br/pquimica/category/revista/Petróleo & Energia, São Paulo, SP, v. 1, n. 1, ; Disponível em: http://www.petroleoenergia.com.br/petroleo/category/revista-petroleo-e-energia/Revista FACTO, Publicação da Associação Brasileira das Indústrias de Química Fina, Biotecnologia e suas Especialidades, Rio de Janeiro, RJ, v. 1, n. 1, ; Disponível em: http://www.abifina.org.br/facto/', 2) | Out-Null
